$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10.82

$ws.Range("D3").Value = 10.19
$ws.Range("E3").Value = 10.8

$ws.Range("B4").Value = 9.18
$ws.Range("C4").Value = 9.81
$ws.Range("E4").Value = 10.63
$ws.Range("F4").Value = 9.72
$ws.Range("J4").Value = 12.67

$ws.Range("C5").Value = 9.15
$ws.Range("D5").Value = 9.37
$ws.Range("F5").Value = 10.14
$ws.Range("H5").Value = 8.37

$ws.Range("D6").Value = 10.28
$ws.Range("E6").Value = 9.86
$ws.Range("H6").Value = 10.47

$ws.Range("H7").Value = 9.88
$ws.Range("J7").Value = 9.53

$ws.Range("E8").Value = 11.63
$ws.Range("F8").Value = 9.53
$ws.Range("G8").Value = 10.12

$ws.Range("D10").Value = 7.33
$ws.Range("G10").Value = 10.47
